# Add a new "Exercises" entry for 2015.11.25 (value 3) as row 3,
# matching the existing Date/Exercises rows above it (A2/B2 = 2015.11.21 / 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 must hold the literal text "2015.11.25" (same text-cell shape as A2),
# not get auto-converted into a date serial number. Stage the text in a
# scratch cell that's formatted as Text first, then copy only the VALUE
# (not the format) into A3 so A3 keeps its plain/default cell style -
# exactly like A2.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "2015.11.25"
$helper.Copy()
$ws.Range("A3").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("B3").Value = 3
